$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- String literals for the new / changed cell contents ---
$row10 = 'Capacitar o aluno a identificar e gerenciar os riscos no ambiente de trabalho, com enfâse nos dos decorrentes das atividades em biotecnologia na industria e em laboratórios de pesquisa.'
$row13 = '8853480 - Tatiane da Franca Silva'
$row14 = 'Introdução a segurança e medicina do trabalho; Conceitos técnico e aspectos legais em biossegurança; Risco em biossegurança ao trabalhador, a comunidade e o ambiente.'
$row16 = 'Gestão em segurança do trabalho e estratégias de prevenção e controle dos riscos no ambiente de trabalho. Normas regulamentadoras;Avaliação de segurança biológica, química e radiológica em atividades de biotecnologia. Manejo e descarte de resíduo contaminado;Critérios e normas regulamentadoras para os diferentes níveis de biossegurança;Legislação para produção e manejo organismos geneticamente modificado e seus derivados;Biossegurança no manuseio de cobaias; Princípios de bioética;Estudos de casos'
$row19 = 'Duas notas N1e N2 distribuídas ao longo do semestre. A composição das "N" fica critério dodocente.'
$row20 = 'MF = (N1 + N2)/2'
$row21 = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor doque 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0.'
$row22 = '1- ABIQUIM . Código de saúde e segurança do trabalhadorguia de implantação.ABIQUIM, 1994. 2- Gonçalves, L. B. Gestão de segurança e medicina do trabalho normas regulamentadoras e fator acidentário de prevenção. Cenofisco, 2011.3- Paoleschi, B. Guia Prático de Segurança do Trabalho. Érica, 2009.4- Bisnfeld, P. C. Biossegurança em Biotecnologia. Interciência, 2004.5- Teixeira. P. e Valle, E. Biossegurança: uma abordagem multidisciplinar, 2002.6- Biosafety in Microbiological and Biomedical Laboratories, 5 ed. U.S. Health Department, 2013.8- Biosecurity, 1ed, 2013.'

# --- Insert a new row at 13 (shifts old rows 13-23 down to 14-24) ---
$ws.Rows("13").Insert()

# The insert carries formatting down into the new row 13 for column A;
# the target layout has no cell in A13 at all, so clear it.
$ws.Range("A13").Clear()

# Give B13/C13 the same number/wrap formatting used by the other B/C body cells
# (style "2" for column B, style "3" for column C) by copying formats down from row 10.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 10 (Objetivos): replace placeholder text with the real objectives text ---
$ws.Range("B10").Value = $row10
$ws.Range("C10").Value = $row10

# --- Row 13 (new row for "Docentes responsaveis" professor name, moved down from row 10) ---
$ws.Range("B13").Value = $row13
$ws.Range("C13").Value = $row13

# --- Row 14 (Programa resumido): replace "Semestral" placeholder with the real summary ---
$ws.Range("B14").Value = $row14
$ws.Range("C14").Value = $row14

# --- Row 16 (Programa): replace the stray date placeholder with the real program text ---
$ws.Range("B16").Value = $row16
$ws.Range("C16").Value = $row16

# --- Row 19 (Metodo): replace placeholder with the real method / criteria text ---
$ws.Range("B19").Value = $row19
$ws.Range("C19").Value = $row19

# --- Row 20 (Criterio): now holds the MF formula text ---
$ws.Range("B20").Value = $row20
$ws.Range("C20").Value = $row20

# --- Row 21 (Norma de recuperacao): now holds the NF formula text; height also shrinks to 60 ---
$ws.Range("B21").Value = $row21
$ws.Range("C21").Value = $row21
$ws.Rows("21").RowHeight = 60

# --- Row 22 (Bibliografia): replace placeholder with the real bibliography text ---
$ws.Range("B22").Value = $row22
$ws.Range("C22").Value = $row22

Write-Output "edit applied"
